$wb = $excel.ActiveWorkbook

# Update OFF sheet, Week row (A2 = "H") with Week 16 data
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 212
$wsOff.Range("C2").Value = 157
$wsOff.Range("D2").Value = 32
$wsOff.Range("E2").Value = 9

# Update DEF sheet, Week row (A2 = "H") with Week 16 data
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 177
$wsDef.Range("C2").Value = 129
$wsDef.Range("D2").Value = 31
$wsDef.Range("E2").Value = 17
$wsDef.Range("F2").Value = 2
